$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy style (bold/border/center) from A147 into the new A column cells
$ws.Range("A147").Copy() | Out-Null
$ws.Range("A148:A152").PasteSpecial(-4122) | Out-Null
# Copy style (datetime number format) from E147 into the new E column cells
$ws.Range("E147").Copy() | Out-Null
$ws.Range("E148:E152").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 148
$ws.Cells.Item(148, 1).Value = 146
$ws.Cells.Item(148, 2).Value = 6788943
$ws.Cells.Item(148, 3).Value = "Croatia HNL"
$ws.Cells.Item(148, 4).Value = "Croatia HNL"
$ws.Cells.Item(148, 5).Value = 45395.58333333334
$ws.Cells.Item(148, 6).Value = "HNK Gorica"
$ws.Cells.Item(148, 7).Value = "Dinamo Zagreb"
$ws.Cells.Item(148, 8).Value = 0
$ws.Cells.Item(148, 9).Value = 2
$ws.Cells.Item(148, 10).Value = "A"
$ws.Cells.Item(148, 11).Value = 8
$ws.Cells.Item(148, 12).Value = 4.5
$ws.Cells.Item(148, 13).Value = 1.363
$ws.Cells.Item(148, 14).Value = 9.5
$ws.Cells.Item(148, 15).Value = 5
$ws.Cells.Item(148, 16).Value = 1.3
$ws.Cells.Item(148, 17).Value = 1.5
$ws.Cells.Item(148, 18).Value = 1.975
$ws.Cells.Item(148, 19).Value = 1.875
$ws.Cells.Item(148, 20).Value = 2.75
$ws.Cells.Item(148, 21).Value = 1.875
$ws.Cells.Item(148, 22).Value = 1.975
$ws.Cells.Item(148, 23).Value = -1
$ws.Cells.Item(148, 24).Value = -1
$ws.Cells.Item(148, 25).Value = 0.3
$ws.Cells.Item(148, 26).Value = -1
$ws.Cells.Item(148, 27).Value = 0.875
$ws.Cells.Item(148, 28).Value = -1
$ws.Cells.Item(148, 29).Value = 0.9750000000000001

# Row 149
$ws.Cells.Item(149, 1).Value = 147
$ws.Cells.Item(149, 2).Value = 6923266
$ws.Cells.Item(149, 3).Value = "Croatia HNL"
$ws.Cells.Item(149, 4).Value = "Croatia HNL"
$ws.Cells.Item(149, 5).Value = 45396.47916666666
$ws.Cells.Item(149, 6).Value = "Hajduk Split"
$ws.Cells.Item(149, 7).Value = "NK Osijek"
$ws.Cells.Item(149, 8).Value = 1
$ws.Cells.Item(149, 9).Value = 2
$ws.Cells.Item(149, 10).Value = "A"
$ws.Cells.Item(149, 11).Value = 1.615
$ws.Cells.Item(149, 12).Value = 3.5
$ws.Cells.Item(149, 13).Value = 6
$ws.Cells.Item(149, 14).Value = 1.85
$ws.Cells.Item(149, 15).Value = 3.2
$ws.Cells.Item(149, 16).Value = 4.5
$ws.Cells.Item(149, 17).Value = -0.5
$ws.Cells.Item(149, 18).Value = 1.85
$ws.Cells.Item(149, 19).Value = 2
$ws.Cells.Item(149, 20).Value = 2.25
$ws.Cells.Item(149, 21).Value = 1.875
$ws.Cells.Item(149, 22).Value = 1.975
$ws.Cells.Item(149, 23).Value = -1
$ws.Cells.Item(149, 24).Value = -1
$ws.Cells.Item(149, 25).Value = 3.5
$ws.Cells.Item(149, 26).Value = -1
$ws.Cells.Item(149, 27).Value = 1
$ws.Cells.Item(149, 28).Value = 0.875
$ws.Cells.Item(149, 29).Value = -1

# Row 150
$ws.Cells.Item(150, 1).Value = 148
$ws.Cells.Item(150, 2).Value = 7993785
$ws.Cells.Item(150, 3).Value = "Croatia HNL"
$ws.Cells.Item(150, 4).Value = "Croatia HNL"
$ws.Cells.Item(150, 5).Value = 45399.5
$ws.Cells.Item(150, 6).Value = "Dinamo Zagreb"
$ws.Cells.Item(150, 7).Value = "NK Varazdin"
$ws.Cells.Item(150, 11).Value = 1.2
$ws.Cells.Item(150, 12).Value = 6.5
$ws.Cells.Item(150, 13).Value = 13
$ws.Cells.Item(150, 14).Value = 1.2
$ws.Cells.Item(150, 15).Value = 6.5
$ws.Cells.Item(150, 16).Value = 13
$ws.Cells.Item(150, 17).Value = -1.75
$ws.Cells.Item(150, 18).Value = 1.8
$ws.Cells.Item(150, 19).Value = 2.05
$ws.Cells.Item(150, 20).Value = 3
$ws.Cells.Item(150, 21).Value = 1.875
$ws.Cells.Item(150, 22).Value = 1.975
$ws.Cells.Item(150, 23).Value = 0
$ws.Cells.Item(150, 24).Value = 0
$ws.Cells.Item(150, 25).Value = 0
$ws.Cells.Item(150, 26).Value = 0
$ws.Cells.Item(150, 27).Value = 0

# Row 151
$ws.Cells.Item(151, 1).Value = 149
$ws.Cells.Item(151, 2).Value = 6954858
$ws.Cells.Item(151, 3).Value = "Croatia HNL"
$ws.Cells.Item(151, 4).Value = "Croatia HNL"
$ws.Cells.Item(151, 5).Value = 45402.59027777778
$ws.Cells.Item(151, 6).Value = "Slaven Belupo"
$ws.Cells.Item(151, 7).Value = "Hajduk Split"
$ws.Cells.Item(151, 11).Value = 5.75
$ws.Cells.Item(151, 12).Value = 4
$ws.Cells.Item(151, 13).Value = 1.571
$ws.Cells.Item(151, 14).Value = 5.5
$ws.Cells.Item(151, 15).Value = 4
$ws.Cells.Item(151, 16).Value = 1.6
$ws.Cells.Item(151, 17).Value = 0.75
$ws.Cells.Item(151, 18).Value = 2.05
$ws.Cells.Item(151, 19).Value = 1.8
$ws.Cells.Item(151, 20).Value = 2.5
$ws.Cells.Item(151, 21).Value = 1.9
$ws.Cells.Item(151, 22).Value = 1.95
$ws.Cells.Item(151, 23).Value = 0
$ws.Cells.Item(151, 24).Value = 0
$ws.Cells.Item(151, 25).Value = 0
$ws.Cells.Item(151, 26).Value = 0
$ws.Cells.Item(151, 27).Value = 0

# Row 152
$ws.Cells.Item(152, 1).Value = 150
$ws.Cells.Item(152, 2).Value = 6962506
$ws.Cells.Item(152, 3).Value = "Croatia HNL"
$ws.Cells.Item(152, 4).Value = "Croatia HNL"
$ws.Cells.Item(152, 5).Value = 45403.59027777778
$ws.Cells.Item(152, 6).Value = "HNK Rijeka"
$ws.Cells.Item(152, 7).Value = "HNK Gorica"
$ws.Cells.Item(152, 11).Value = 1.285
$ws.Cells.Item(152, 12).Value = 6
$ws.Cells.Item(152, 13).Value = 8.5
$ws.Cells.Item(152, 14).Value = 1.25
$ws.Cells.Item(152, 15).Value = 6
$ws.Cells.Item(152, 16).Value = 9
$ws.Cells.Item(152, 17).Value = -1.75
$ws.Cells.Item(152, 18).Value = 2.025
$ws.Cells.Item(152, 19).Value = 1.825
$ws.Cells.Item(152, 20).Value = 3
$ws.Cells.Item(152, 21).Value = 2
$ws.Cells.Item(152, 22).Value = 1.85
$ws.Cells.Item(152, 23).Value = 0
$ws.Cells.Item(152, 24).Value = 0
$ws.Cells.Item(152, 25).Value = 0
$ws.Cells.Item(152, 26).Value = 0
$ws.Cells.Item(152, 27).Value = 0
